# "Root and Blind are final"
# Populates the "End Effect + Notes" (L) column with "trivial" for every
# section except the auto-attack-in-air interaction row (L13), which gets
# the full explanatory note. Also adds the "final notes" block under the
# table (rows 70-72) and merges C70:G70 for the first note line.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- sheet view: keep selection near the newly-added notes section -------
$ws.Range("C74").Select()
$excel.ActiveWindow.ScrollRow = 61
$excel.ActiveWindow.ScrollColumn = 1

# --- narrow spacer columns -------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 0.79
$ws.Columns.Item(16).ColumnWidth = 55.94

# --- "trivial" for every CC-interaction row except the auto-attack one ---
$trivialRows = @(
    11,12,   14,15,16,17,18,
    21,22,23,24,25,26,27,28,
    31,32,33,34,35,36,37,38,
    41,42,43,44,45,46,47,48,
    51,52,53,54,55,56,57,58,
    61,62,63,64,65,66,67,68
)
foreach ($r in $trivialRows) {
    $ws.Cells.Item($r, 12).Value = "trivial"
}

# Auto attack in range + cc: the one row with real behavioural nuance.
$ws.Cells.Item(13, 12).Value = "if blinded while auto is in air, its considered a miss. If the debuff times out and an auto is in the air, it still misses"

# --- final notes block (rows 70-72) ---------------------------------------
$ws.Cells.Item(70, 2).Value = "final notes: "
$ws.Cells.Item(71, 3).Value = "If the debuff times out and an auto is in the air, it still misses"
$ws.Cells.Item(70, 3).Value = "If blinded while auto is in air, its considered a miss"
$ws.Cells.Item(72, 3).Value = "Everything else has no differences from the normal behavior"

$notesRange = $ws.Range("C70:G70")
$notesRange.HorizontalAlignment = -4131
$notesRange.VerticalAlignment = -4108
$notesRange.Merge()

Write-Output "edit complete"
